$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: was CNN_128_filters_3_kernels_predictions / correction -> now mirrors row 2's
# CNN_256_filters_3_kernels_predictions / correction data, label "Best Recall, Best F1"
$ws.Range("A4").Value = "CNN_256_filters_3_kernels_predictions"
$ws.Range("B4").Value = "correction"
$ws.Range("D4").Value = 0.8225806451612904
$ws.Range("E4").Value = 0.6538461538461539
$ws.Range("F4").Value = 0.7285714285714286
$ws.Range("G4").Value = 0.8225806451612904
$ws.Range("H4").Value = 0.6538461538461539
$ws.Range("I4").Value = 0.7285714285714286
$ws.Range("K4").Value = 0.8095238095238095
$ws.Range("L4").Value = 0.6538461538461539
$ws.Range("M4").Value = 0.7234042553191489
$ws.Range("N4").Value = -0.01305683563748083
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0.005167173252279778
$ws.Range("Q4").Value = -0.01587301587301592
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = -0.007092198581560478
$ws.Range("T4").Value = "Best Recall, Best F1"

# Row 5: was CNN_128_filters_3_kernels_predictions / detection_correction -> now mirrors row 3's
# CNN_256_filters_3_kernels_predictions / detection_correction data, label "Best Recall, Best F1"
$ws.Range("A5").Value = "CNN_256_filters_3_kernels_predictions"
$ws.Range("B5").Value = "detection_correction"
$ws.Range("D5").Value = 0.8225806451612904
$ws.Range("E5").Value = 0.6538461538461539
$ws.Range("F5").Value = 0.7285714285714286
$ws.Range("G5").Value = 0.8225806451612904
$ws.Range("H5").Value = 0.6538461538461539
$ws.Range("I5").Value = 0.7285714285714286
$ws.Range("K5").Value = 0.8472222222222222
$ws.Range("L5").Value = 0.782051282051282
$ws.Range("M5").Value = 0.8133333333333332
$ws.Range("N5").Value = 0.02464157706093184
$ws.Range("O5").Value = 0.1282051282051282
$ws.Range("P5").Value = 0.08476190476190459
$ws.Range("Q5").Value = 0.02995642701525047
$ws.Range("R5").Value = 0.196078431372549
$ws.Range("S5").Value = 0.1163398692810455
$ws.Range("T5").Value = "Best Recall, Best F1"

# Row 6: was CNN_256_filters_5_kernels_predictions / correction -> now
# CNN_Attention_128_filters_5_kernels_predictions / correction with new values + new T label
$ws.Range("A6").Value = "CNN_Attention_128_filters_5_kernels_predictions"
$ws.Range("B6").Value = "correction"
$ws.Range("D6").Value = 0.9375
$ws.Range("E6").Value = 0.1923076923076923
$ws.Range("F6").Value = 0.3191489361702128
$ws.Range("G6").Value = 0.9375
$ws.Range("H6").Value = 0.1923076923076923
$ws.Range("I6").Value = 0.3191489361702128
$ws.Range("K6").Value = 0.9375
$ws.Range("L6").Value = 0.1923076923076923
$ws.Range("M6").Value = 0.3191489361702128
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 2.775557561562891 / 100000000000000000
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 1.443289932012704 / 10000000000000000
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = "Best Precision"

# Row 7: was CNN_256_filters_5_kernels_predictions / detection_correction -> now
# CNN_Attention_128_filters_5_kernels_predictions / detection_correction with new values + new T label
$ws.Range("A7").Value = "CNN_Attention_128_filters_5_kernels_predictions"
$ws.Range("B7").Value = "detection_correction"
$ws.Range("D7").Value = 0.9375
$ws.Range("E7").Value = 0.1923076923076923
$ws.Range("F7").Value = 0.3191489361702128
$ws.Range("G7").Value = 0.9375
$ws.Range("H7").Value = 0.1923076923076923
$ws.Range("I7").Value = 0.3191489361702128
$ws.Range("K7").Value = 0.9375
$ws.Range("L7").Value = 0.1923076923076923
$ws.Range("M7").Value = 0.3191489361702128
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 2.775557561562891 / 100000000000000000
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 1.443289932012704 / 10000000000000000
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = "Best Precision"
